$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131, pushing the existing rows 131-136 down to 132-137.
$ws.Rows(131).Insert()

# Populate the newly inserted row 131 with the weekly price record
# (same Feria/region/product as its neighbours, new date 2021-11-16 / 44516).
$ws.Range("A131").Value = 4
$ws.Range("B131").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C131").Value = "Los Lagos"
$ws.Range("D131").Value = 44516
$ws.Range("E131").Value = 10
$ws.Range("F131").Value = 100112039
$ws.Range("G131").Value = "Ciboulette"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 240
$ws.Range("K131").Value = 2500
$ws.Range("L131").Value = 2500
$ws.Range("M131").Value = 2500
$ws.Range("N131").Value = "$/docena de atados"
$ws.Range("O131").Value = "Región Metropolitana"
$ws.Range("P131").Value = 833
$ws.Range("Q131").Value = 3
$ws.Range("R131").Value = "Hortaliza"
